{"js": "// Update the two-digit x two-digit multiplication answers in the table.\nconst replacements = [\n  [\"50\u00d799=4950\", \"45\u00d736=1620\"],\n  [\"86\u00d759=5074\", \"79\u00d725=1975\"],\n  [\"46\u00d711=506\", \"57\u00d760=3420\"],\n  [\"23\u00d770=1610\", \"99\u00d769=6831\"],\n  [\"36\u00d749=1764\", \"79\u00d716=1264\"],\n  [\"23\u00d789=2047\", \"47\u00d736=1692\"],\n  [\"37\u00d749=1813\", \"43\u00d793=3999\"],\n  [\"28\u00d720=560\", \"14\u00d713=182\"],\n  [\"88\u00d716=1408\", \"63\u00d754=3402\"],\n  [\"40\u00d742=1680\", \"14\u00d782=1148\"],\n  [\"36\u00d784=3024\", \"61\u00d790=5490\"],\n  [\"40\u00d725=1000\", \"85\u00d767=5695\"],\n  [\"36\u00d793=3348\", \"80\u00d730=2400\"],\n  [\"40\u00d758=2320\", \"41\u00d797=3977\"],\n  [\"98\u00d787=8526\", \"79\u00d766=5214\"],\n  [\"76\u00d762=4712\", \"17\u00d797=1649\"],\n  [\"32\u00d772=2304\", \"64\u00d752=3328\"],\n  [\"95\u00d758=5510\", \"59\u00d783=4897\"],\n  [\"50\u00d783=4150\", \"84\u00d726=2184\"],\n  [\"87\u00d796=8352\", \"20\u00d753=1060\"],\n  [\"85\u00d754=4590\", \"62\u00d726=1612\"],\n  [\"83\u00d785=7055\", \"94\u00d765=6110\"],\n  [\"79\u00d753=4187\", \"80\u00d751=4080\"],\n  [\"90\u00d773=6570\", \"75\u00d713=975\"],\n  [\"68\u00d781=5508\", \"99\u00d730=2970\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"50\u00d799=4950\", \"45\u00d736=1620\"),\n    @(\"86\u00d759=5074\", \"79\u00d725=1975\"),\n    @(\"46\u00d711=506\", \"57\u00d760=3420\"),\n    @(\"23\u00d770=1610\", \"99\u00d769=6831\"),\n    @(\"36\u00d749=1764\", \"79\u00d716=1264\"),\n    @(\"23\u00d789=2047\", \"47\u00d736=1692\"),\n    @(\"37\u00d749=1813\", \"43\u00d793=3999\"),\n    @(\"28\u00d720=560\", \"14\u00d713=182\"),\n    @(\"88\u00d716=1408\", \"63\u00d754=3402\"),\n    @(\"40\u00d742=1680\", \"14\u00d782=1148\"),\n    @(\"36\u00d784=3024\", \"61\u00d790=5490\"),\n    @(\"40\u00d725=1000\", \"85\u00d767=5695\"),\n    @(\"36\u00d793=3348\", \"80\u00d730=2400\"),\n    @(\"40\u00d758=2320\", \"41\u00d797=3977\"),\n    @(\"98\u00d787=8526\", \"79\u00d766=5214\"),\n    @(\"76\u00d762=4712\", \"17\u00d797=1649\"),\n    @(\"32\u00d772=2304\", \"64\u00d752=3328\"),\n    @(\"95\u00d758=5510\", \"59\u00d783=4897\"),\n    @(\"50\u00d783=4150\", \"84\u00d726=2184\"),\n    @(\"87\u00d796=8352\", \"20\u00d753=1060\"),\n    @(\"85\u00d754=4590\", \"62\u00d726=1612\"),\n    @(\"83\u00d785=7055\", \"94\u00d765=6110\"),\n    @(\"79\u00d753=4187\", \"80\u00d751=4080\"),\n    @(\"90\u00d773=6570\", \"75\u00d713=975\"),\n    @(\"68\u00d781=5508\", \"99\u00d730=2970\"),\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair[0]\n    $find.Replacement.Text = $pair[1]\n    $find.Execute($null, $true, $true, $false, $false, $false, $true, 1, $false, $pair[1], 2) | Out-Null\n}\n"}
